$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.801.48'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = '2.345.75'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.97'
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("E6").Value = '  -4.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.35'
$ws.Range("E7").Value = '  -3.58%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -3.04%  '

$ws.Range("E10").Value = '  -0.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.26'
$ws.Range("E11").Value = '  +4.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.91'
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.24'
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").Value = '2.695.65'
$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.07'
$ws.Range("E16").Value = '  -3.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.903'
$ws.Range("E17").Value = '  -2.21%  '

$ws.Range("D18").Value = '2.347.78'
$ws.Range("E18").Value = '  +0.52%  '

$ws.Range("D19").Value = '43.759.82'
$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("E20").Value = '  +1.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '78.34'
$ws.Range("E21").Value = '  +1.06%  '

$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.66'
$ws.Range("E23").Value = '  -2.42%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.78'
$ws.Range("E25").Value = '  +1.16%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.50'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.39'
$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  +0.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.34'
$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.20'
$ws.Range("E31").Value = '  -4.06%  '

$ws.Range("E32").Value = '  -0.08%  '

$ws.Range("E33").Value = '  -2.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0746'
$ws.Range("E34").Value = '  -1.99%  '

$ws.Range("E35").Value = '  -5.86%  '

$ws.Range("E36").Value = '  -1.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("E37").Value = '  +1.54%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.37'
$ws.Range("E39").Value = '  -0.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.61'
$ws.Range("E40").Value = '  +14.84%  '

$ws.Range("E41").Value = '  -4.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.87'
$ws.Range("E42").Value = '  +16.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.17'
$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.74'
$ws.Range("E44").Value = '  -2.33%  '

$ws.Range("E45").Value = '  -3.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.197'
$ws.Range("E46").Value = '  -9.88%  '

$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("E48").Value = '  -2.75%  '

$ws.Range("E49").Value = '  -3.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.15'
$ws.Range("E50").Value = '  -3.84%  '

$ws.Range("E51").Value = '  -5.33%  '
